$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ciudades")

# Update the "last updated" timestamp text in A1
$ws.Range("A1").Value = "Datos actualizados a 9 de Abril de 2020 a las 23:22"

# Update the Cataluña row (row 5) figures
$ws.Range("B5").Value = 31727
$ws.Range("C5").Value = 13513
$ws.Range("D5").Value = 14983
$ws.Range("E5").Value = 3231
